$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated query text for the "ParticipantsTab" row (B2) ---
# Adds a `with p` line before the OPTIONAL MATCH on study, drops the
# trailing space after RETURN, and upper-cases the ORDER BY / LIMIT clause.
$participantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE COALESCE(g.platform, "Not specified in data") in ['Not specified in data']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

# --- Updated query text for the "FilesTab" row (B4) ---
# Adds an extra space after WHERE.
$fileQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE  COALESCE(g.platform, "Not specified in data") in ['Not specified in data']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name limit 100
'@

# --- New query text for the "SamplesTab" row (B3) ---
# Replaces the old Sample ID query entirely with a rewritten version that
# matches from sample outward, filters distinct (s,p,samp) and returns the
# Sample ID first.
$sampleQuery = @'
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE COALESCE(g.platform, "Not specified in data") in ['Not specified in data']
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(samp.sample_tumor_status,'') as `Tumor`,
    coalesce(samp.sample_type,'') as `Analyte Type`
ORDER BY samp.sample_id limit 100
'@

# Row 2 = ParticipantsTab, Row 3 = SamplesTab, Row 4 = FilesTab (col A).
# Column B holds the tab-specific query, column C holds the shared
# summary-stats query (unchanged, just shifts shared-string index on save).
# Write the file-query update before the sample-query replacement so the
# brand-new sample-ID string lands last in the shared-string table (matching
# the upstream diff's append order).
$ws.Range("B2").Value = $participantQuery
$ws.Range("B4").Value = $fileQuery
$ws.Range("B3").Value = $sampleQuery

# --- View state: scroll so row 2 is at the top and select C2 ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("C2").Select()
